$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the "Actief" (Active) flag for the two groups:
# - "Totaal" group becomes active
# - "test" group becomes inactive (can now be considered removed/deleted)
$ws.Range("I2").Value = $true
$ws.Range("I3").Value = $false

# Remove trailing empty rows that used to reserve space for extra groups.
$ws.Rows("4:8").Delete()
